# "All models, plus correlation"
# Rename the two worksheets, add six more model rows (and matching
# classifier names) to the Results sheet, widen column B to fit the
# longer names, and add a third Feature/Importance block ("Model 9")
# to the Importance sheet for the Balanced Random Forest Classifier.

$wb = $excel.ActiveWorkbook

# --- Rename the worksheets ------------------------------------------------
$wsResults = $wb.Worksheets.Item("Results")
$wsResults.Name = "Red Results"

$wsImportance = $wb.Worksheets.Item("Importance")
$wsImportance.Name = "Red Importance"

# --- Results sheet: add the six additional models -------------------------
$results = @(
    @(5,  "Gradient Boosted Tree",              0.89, 0.56, 0.43, 0.49),
    @(6,  "Logisitic-ROS",                      0.83, 0.35, 0.9,  0.51),
    @(7,  "Logistic-SMOTE",                     0.84, 0.37, 0.92, 0.53),
    @(8,  "Logistic-Cluster centroid US",       0.82, 0.34, 0.88, 0.49),
    @(9,  "Balanced Random Forest Classifier",  0.84, 0.36, 0.92, 0.52),
    @(10, "Easy Ensemble AdaBoos Classifier",   0.82, 0.33, 0.9,  0.48)
)

$row = 6
foreach ($r in $results) {
    $wsResults.Cells.Item($row, 1).Value = $r[0]
    $wsResults.Cells.Item($row, 2).Value = $r[1]
    $wsResults.Cells.Item($row, 3).Value = $r[2]
    $wsResults.Cells.Item($row, 4).Value = $r[3]
    $wsResults.Cells.Item($row, 5).Value = $r[4]
    $wsResults.Cells.Item($row, 6).Value = $r[5]
    $row++
}

$wsResults.Columns.Item(2).ColumnWidth = 31.16666666666666
$wsResults.Range("C16").Select() | Out-Null

# --- Importance sheet: add the "Model 9" feature-importance block ---------
$wsImportance.Cells.Item(1, 10).Value = "Model 9"
$wsImportance.Cells.Item(2, 10).Value = "Feature"
$wsImportance.Cells.Item(2, 11).Value = "Importance"

$importance = @(
    @("ph",                    0.233499),
    @("volatile_acidity",      0.128257),
    @("alcohol",                0.116748),
    @("citric_acid",            0.083079),
    @("density",                0.076444),
    @("total_sulfur_dioxide",   0.074452),
    @("fixed_acidity",          0.070034),
    @("chlorides",               0.067913),
    @("sulphates",               0.053655),
    @("residual_sugar",          0.049013),
    @("free_sulfur_dioxide",     0.046905)
)

$row = 3
foreach ($i in $importance) {
    $wsImportance.Cells.Item($row, 10).Value = $i[0]
    $wsImportance.Cells.Item($row, 11).Value = $i[1]
    $row++
}

$wsImportance.Columns.Item(10).ColumnWidth = 18.5
$wsImportance.Columns.Item(11).ColumnWidth = 10.33
